# Insert a new data row at row 106 (a new weekly price observation),
# pushing the existing rows 106-159 down to 107-160.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(106).Insert()

# Populate the newly inserted row with the latest observation.
$ws.Range("A106").Value = 5
$ws.Range("B106").Value = "Macroferia Regional de Talca"
$ws.Range("C106").Value = "Maule"
$ws.Range("D106").Value = 44719
$ws.Range("E106").Value = 7
$ws.Range("F106").Value = 100112031
$ws.Range("G106").Value = "Poroto verde"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 150
$ws.Range("K106").Value = 23000
$ws.Range("L106").Value = 23000
$ws.Range("M106").Value = 23000
$ws.Range("N106").Value = "$/malla 25 kilos"
$ws.Range("O106").Value = "Región de Arica y Parinacota"
$ws.Range("P106").Value = 920
$ws.Range("Q106").Value = 25
$ws.Range("R106").Value = "Hortaliza"
